$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. "Handoff transform failed" -> "Ready for handoff"
#    This text is a shared string used by Overview!B2, Overview!C2,
#    zh-cn!B2 and de-de!B2 -- update every cell that carries it so the
#    underlying shared string itself changes (report status is now ready).
# ---------------------------------------------------------------------------
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("B2").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: report the generated handoff artifact for the zh-cn target
# ---------------------------------------------------------------------------
$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/de9d84eb10155872ade5411a9c6aac7de111129e"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "$repoBase/0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.3c01fdf2275413683ef44bdf663d5525faa94af6.zh-cn.xlf", "", "", "0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.3c01fdf2275413683ef44bdf663d5525faa94af6.zh-cn.xlf")
$wsZhCn.Range("D2").Value = "2016-01-11 14:11:51"
$wsZhCn.Range("H2").Value = "Include"

# ---------------------------------------------------------------------------
# 3. de-de sheet: report the generated handoff artifact for the de-de target
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "$repoBase/0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.3c01fdf2275413683ef44bdf663d5525faa94af6.de-de.xlf", "", "", "0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.3c01fdf2275413683ef44bdf663d5525faa94af6.de-de.xlf")
$wsDeDe.Range("D2").Value = "2016-01-11 14:12:10"
$wsDeDe.Range("H2").Value = "Include"
